# Actualiza base de datos EC: reordena los periodos de mora de la columna
# "Periodo Mora" (E16:E20) de orden descendente a orden ascendente.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2011"
$ws.Range("E17").Value = "2012"
$ws.Range("E18").Value = "2101"
$ws.Range("E19").Value = "2102"
$ws.Range("E20").Value = "2103"
